# authorization.xlsx — refactor "true_false" choice list into "true_false_unsure"
# and add the new "Unsure" choice, per commit:
#   "Refactor logic to reduce duplicate code, include unsure option as
#    verification/authorization choice"

$wb = $excel.ActiveWorkbook

$choices = $wb.Worksheets.Item("choices")
$survey  = $wb.Worksheets.Item("survey")

# --- choices sheet -------------------------------------------------------
# Rename the existing "true_false" choice_list to "true_false_unsure" for
# both of its current rows (False / True).
$choices.Range("A6").Value = "true_false_unsure"
$choices.Range("A7").Value = "true_false_unsure"

# Insert a new row right after "True" (row 7) for the new "Unsure" choice,
# reusing the blank spacer row's slot so the male_female table below shifts
# down naturally.
$choices.Rows.Item(8).Insert()
$choices.Range("A8").Value = "true_false_unsure"
$choices.Range("B8").Value = -1
$choices.Range("C8").Value = "Unsure"

# B8 keeps the plain/default number format (unlike B6/B7's integer format).
$choices.Range("B8").ClearFormats()

# The blank spacer row's height (19) now belongs to the new data row; the
# still-blank row right below it (row 9) keeps that same height too.
$choices.Rows.Item(8).RowHeight = 19

# --- survey sheet ----------------------------------------------------------
# The "Has the VEO authorized this business?" question now points at the
# renamed choice list.
$survey.Range("E5").Value = "true_false_unsure"

# --- cosmetic view state (best effort) -------------------------------------
$choices.Activate()
$choices.Range("C13").Select()

try {
    $win = $excel.ActiveWindow
    $win.Left = 5320
    $win.Top = 760
} catch {
}

Write-Host "applied true_false -> true_false_unsure refactor"
